$d = $word.ActiveDocument

# New bullet points to insert right after the "Full-Stack Development and
# Data Engineering" sub-heading under the Siege Analytics (PARTNER) entry,
# and before the existing "Lead comprehensive research initiatives..."
# bullet.
$newBullets = @(
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times",
    "• Architected systems supporting 2,500+ concurrent users conducting redistricting analysis",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

# Locate the anchor paragraph ("Full-Stack Development and Data Engineering")
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Full-Stack Development and Data Engineering") {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    foreach ($t in $newBullets) {
        $anchor.Range.InsertParagraphAfter()
        $anchor = $anchor.Next()
        $anchor.Range.Text = $t
    }
}
